# C5-PowerPoint.pptx — Fri, Apr 03, 2020 7:06:34 PM
#
# 1) Slide 6: table picks up a different (built-in) table style.
# 2) The deck's theme (currently the custom "Integral" palette) is switched
#    back to the default "Office" colour scheme via the Design/Theme colours.

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on slide 6 -------------------------------------
$tableSlide = $p.Slides.Item(6)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shp = $tableSlide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{FB61310B-BF97-4BF4-BF3E-24F480C16070}")
    }
}

# --- 2. Switch the theme colours from "Integral" back to "Office" ---------
function RGBVal($r, $g, $b) { return $r + ($g * 256) + ($b * 65536) }

# Office theme colour scheme, in ppColorSchemeIndex order:
#   1 Background1(dk1) 2 Text1(lt1) 3 Background2(dk2) 4 Text2(lt2)
#   5-10 Accent1-6      11 Hyperlink 12 FollowedHyperlink
$officeColors = @(
    (RGBVal 0x00 0x00 0x00),  # dk1
    (RGBVal 0xFF 0xFF 0xFF),  # lt1
    (RGBVal 0x44 0x54 0x6A),  # dk2
    (RGBVal 0xE7 0xE6 0xE6),  # lt2
    (RGBVal 0x5B 0x9B 0xD5),  # accent1
    (RGBVal 0xED 0x7D 0x31),  # accent2
    (RGBVal 0xA5 0xA5 0xA5),  # accent3
    (RGBVal 0xFF 0xC0 0x00),  # accent4
    (RGBVal 0x44 0x72 0xC4),  # accent5
    (RGBVal 0x70 0xAD 0x47),  # accent6
    (RGBVal 0x05 0x63 0xC1),  # hyperlink
    (RGBVal 0x95 0x4F 0x72)   # followed hyperlink
)

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Item($i).RGB = $officeColors[$i - 1]
}
